$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 741.1875
$ws.Range("I8").Value = 96.84614999999999
$ws.Range("J8").Value = 3533.3333
$ws.Range("K8").Value = 290.53845
$ws.Range("L8").Value = 10599.9999
$ws.Range("M8").Value = -151.53845
$ws.Range("N8").Value = -10877.9999

$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = ""

$ws.Range("H13").Value = 29750
$ws.Range("J13").Value = 29750
$ws.Range("L13").Value = 29750
$ws.Range("N13").Value = -30088

$ws.Range("H28").Value = 1541
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = ""

$ws.Range("H64").Value = 3271.4285
$ws.Range("I64").Value = 2502
$ws.Range("J64").Value = 3399.6667
$ws.Range("K64").Value = 2502
$ws.Range("L64").Value = 3399.6667
$ws.Range("M64").Value = -2254
$ws.Range("N64").Value = -3895.6667

$ws.Range("H67").Value = 3271.4285
$ws.Range("I67").Value = 2502
$ws.Range("J67").Value = 3399.6667
$ws.Range("K67").Value = 2502
$ws.Range("L67").Value = 3399.6667
$ws.Range("M67").Value = -1644
$ws.Range("N67").Value = -5115.6667

$ws.Range("H113").Value = 5550.4
$ws.Range("J113").Value = 5550.4
$ws.Range("L113").Value = 5550.4
$ws.Range("N113").Value = -12058.4

$ws.Range("H138").Value = 2193.61
$ws.Range("I138").Value = 556.9655
$ws.Range("J138").Value = 2862.0986
$ws.Range("K138").Value = 1670.8965
$ws.Range("L138").Value = 8586.2958
$ws.Range("M138").Value = 3469.1035
$ws.Range("N138").Value = -18866.2958

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1459.8
$ws.Range("I2").Value = 1181.6
$ws.Range("J2").Value = 1738
$ws.Range("K2").Value = 1181.6
$ws.Range("L2").Value = 1738
$ws.Range("M2").Value = -1068.6
$ws.Range("N2").Value = -1964

$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = ""

$ws.Range("H10").Value = 15925.333
$ws.Range("J10").Value = 15925.333
$ws.Range("L10").Value = 15925.333
$ws.Range("N10").Value = -16265.333

$ws.Range("H32").Value = 6523.7046
$ws.Range("I32").Value = 5610.921
$ws.Range("J32").Value = 12304.667
$ws.Range("K32").Value = 5610.921
$ws.Range("L32").Value = 12304.667
$ws.Range("M32").Value = -5323.921
$ws.Range("N32").Value = -12878.667

$ws.Range("H61").Value = 1505.3636
$ws.Range("I61").Value = 1258.875
$ws.Range("K61").Value = 1258.875
$ws.Range("M61").Value = -1046.875

$ws.Range("H116").Value = 1459.8
$ws.Range("I116").Value = 1181.6
$ws.Range("J116").Value = 1738
$ws.Range("K116").Value = 1181.6
$ws.Range("L116").Value = 1738
$ws.Range("M116").Value = 1112.4
$ws.Range("N116").Value = -6326

$ws.Range("H136").Value = 1505.3636
$ws.Range("I136").Value = 1258.875
$ws.Range("K136").Value = 3776.625
$ws.Range("M136").Value = -1226.625

$ws.Range("H137").Value = 44560
$ws.Range("J137").Value = 44560
$ws.Range("L137").Value = 44560
$ws.Range("N137").Value = -54760

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1459.8
$ws.Range("I3").Value = 1181.6
$ws.Range("J3").Value = 1738
$ws.Range("K3").Value = 1181.6
$ws.Range("L3").Value = 1738
$ws.Range("M3").Value = -1067.6
$ws.Range("N3").Value = -1966

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2857.4285
$ws.Range("I4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("M4").Value = 111

$ws.Range("H11").Value = 26333.334
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 26333.334
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 26333.334
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = -26613.334

$ws.Range("H31").Value = 4261.1665
$ws.Range("I31").Value = 1283.5
$ws.Range("J31").Value = 5750
$ws.Range("K31").Value = 1283.5
$ws.Range("L31").Value = 5750
$ws.Range("M31").Value = -988.5
$ws.Range("N31").Value = -6340

$ws.Range("H34").Value = 4261.1665
$ws.Range("I34").Value = 1283.5
$ws.Range("J34").Value = 5750
$ws.Range("K34").Value = 1283.5
$ws.Range("L34").Value = 5750
$ws.Range("M34").Value = -1081.5
$ws.Range("N34").Value = -6154

$ws.Range("H39").Value = 21460
$ws.Range("I39").Value = 2500
$ws.Range("K39").Value = 2500
$ws.Range("M39").Value = -2109

$ws.Range("H49").Value = 21460
$ws.Range("I49").Value = 2500
$ws.Range("K49").Value = 2500
$ws.Range("M49").Value = -2318

$ws.Range("H86").Value = 2519.25
$ws.Range("I86").Value = 2136.4546
$ws.Range("J86").Value = 3361.4
$ws.Range("K86").Value = 2136.4546
$ws.Range("L86").Value = 3361.4
$ws.Range("M86").Value = -1013.4546
$ws.Range("N86").Value = -5607.4

$ws.Range("H89").Value = 2519.25
$ws.Range("I89").Value = 2136.4546
$ws.Range("J89").Value = 3361.4
$ws.Range("K89").Value = 10682.273
$ws.Range("L89").Value = 16807
$ws.Range("M89").Value = -5066.273000000001
$ws.Range("N89").Value = -28039

$ws.Range("H94").Value = 877.5
$ws.Range("I94").Value = 489.53845
$ws.Range("J94").Value = 1265.4615
$ws.Range("K94").Value = 489.53845
$ws.Range("L94").Value = 1265.4615
$ws.Range("M94").Value = -38.53845000000001
$ws.Range("N94").Value = -2167.4615

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 7993.1665
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 9391.799999999999
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 28175.4
$ws.Range("M69").Value = -2189
$ws.Range("N69").Value = -29797.4

$ws.Range("H72").Value = 7993.1665
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 9391.799999999999
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 84526.2
$ws.Range("M72").Value = -4944
$ws.Range("N72").Value = -92638.2

$ws.Range("H97").Value = 503.1111
$ws.Range("I97").Value = 130
$ws.Range("J97").Value = 689.6667
$ws.Range("K97").Value = 390
$ws.Range("L97").Value = 2069.0001
$ws.Range("M97").Value = 106
$ws.Range("N97").Value = -3061.0001

$ws.Range("H113").Value = 5000846.5
$ws.Range("I113").Value = 689.5
$ws.Range("J113").Value = 9616376
$ws.Range("K113").Value = 2068.5
$ws.Range("L113").Value = 28849128
$ws.Range("M113").Value = 101.5
$ws.Range("N113").Value = -28853468

$ws.Range("H131").Value = 791.8200000000001
$ws.Range("I131").Value = 509.66666
$ws.Range("J131").Value = 800.5463999999999
$ws.Range("K131").Value = 1528.99998
$ws.Range("L131").Value = 2401.6392
$ws.Range("M131").Value = 3511.00002
$ws.Range("N131").Value = -12481.6392

$ws.Range("H132").Value = 2288.7778
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 5639.6
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 50756.4
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -55816.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = ""

$ws.Range("H14").Value = 10926838
$ws.Range("I14").Value = 12909282
$ws.Range("J14").Value = 23399.5
$ws.Range("K14").Value = 12909282
$ws.Range("L14").Value = 23399.5
$ws.Range("M14").Value = -12909114
$ws.Range("N14").Value = -23735.5

$ws.Range("H80").Value = 50002624
$ws.Range("I80").Value = 83335370
$ws.Range("J80").Value = 3503
$ws.Range("K80").Value = 83335370
$ws.Range("L80").Value = 3503
$ws.Range("M80").Value = -83334372
$ws.Range("N80").Value = -5499

$ws.Range("H83").Value = 50002624
$ws.Range("I83").Value = 83335370
$ws.Range("J83").Value = 3503
$ws.Range("K83").Value = 416676850
$ws.Range("L83").Value = 17515
$ws.Range("M83").Value = -416671858
$ws.Range("N83").Value = -27499

$ws.Range("H105").Value = 32000
$ws.Range("J105").Value = 32000
$ws.Range("L105").Value = 32000
$ws.Range("N105").Value = -38988

$ws.Range("H137").Value = 38643.332
$ws.Range("J137").Value = 38643.332
$ws.Range("L137").Value = 38643.332
$ws.Range("N137").Value = -48843.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 51522.2
$ws.Range("I23").Value = 39200
$ws.Range("K23").Value = 39200
$ws.Range("M23").Value = -38971

$ws.Range("H113").Value = 352
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 352
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1056
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -5396
